$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 24.19000000000034
$ws.Range("H2").Value = 0.000490479086283413
$ws.Range("I2").Value = 0.000490479086283413
$ws.Range("L2").Value = 47.5575205468682
$ws.Range("M2").Value = "[19.001831568798224, 76.11320952493818]"
$ws.Range("N2").Value = 0.001622621842740646
$ws.Range("O2").Value = 0.001622621842740646
$ws.Range("P2").Value = 1.805079262422964
$ws.Range("Q2").Value = "[1.0880791372793484, 2.5220793875665803]"
$ws.Range("R2").Value = [double]"7.269495455375363e-06"
$ws.Range("S2").Value = [double]"7.269495455375363e-06"
$ws.Range("T2").Value = 59.00631883470313
$ws.Range("U2").Value = "[43.04926774208796, 74.9633699273183]"
$ws.Range("V2").Value = [double]"2.219213035559164e-09"
$ws.Range("W2").Value = [double]"2.219213035559164e-09"
$ws.Range("X2").Value = 17.24052052052076
$ws.Range("Y2").Value = 14.48010010010031
$ws.Range("Z2").Value = 20.00094094094122

# Row 3
$ws.Range("F3").Value = 24.19000000000034
$ws.Range("H3").Value = 0.0001105256337661231
$ws.Range("I3").Value = 0.0001105256337661231
$ws.Range("L3").Value = 44.59693068219404
$ws.Range("M3").Value = "[20.263405493064667, 68.9304558713234]"
$ws.Range("N3").Value = 0.0006004799787404647
$ws.Range("O3").Value = 0.0006004799787404647
$ws.Range("P3").Value = 1.213868632918579
$ws.Range("Q3").Value = "[0.5975001042863468, 1.8302371615508104]"
$ws.Range("R3").Value = 0.0002588134454946101
$ws.Range("S3").Value = 0.0002588134454946101
$ws.Range("T3").Value = 57.36105963487238
$ws.Range("U3").Value = "[43.935059142484896, 70.78706012725988]"
$ws.Range("V3").Value = [double]"4.615641202576626e-11"
$ws.Range("W3").Value = [double]"4.615641202576626e-11"
$ws.Range("X3").Value = 19.51665665665693
$ws.Range("Y3").Value = 17.14366366366391
$ws.Range("Z3").Value = 21.88964964964996

# Row 4
$ws.Range("F4").Value = 24.19000000000034
$ws.Range("H4").Value = 0.03471946667985881
$ws.Range("I4").Value = 0.03471946667985881
$ws.Range("L4").Value = 23.55730984922108
$ws.Range("M4").Value = "[1.4644946803476273, 45.65012501809453]"
$ws.Range("N4").Value = 0.03716639992294191
$ws.Range("O4").Value = 0.03716639992294191
$ws.Range("P4").Value = 1.352237078121733
$ws.Range("Q4").Value = "[0.03144737390980801, 2.6730267823336575]"
$ws.Range("R4").Value = 0.04500473056051235
$ws.Range("S4").Value = 0.04500473056051235
$ws.Range("T4").Value = 48.80677474593714
$ws.Range("U4").Value = "[36.30716619987752, 61.30638329199676]"
$ws.Range("V4").Value = [double]"5.433984373581779e-10"
$ws.Range("W4").Value = [double]"5.433984373581779e-10"
$ws.Range("X4").Value = 18.98394394394421
$ws.Range("Y4").Value = 13.89895895895916
$ws.Range("Z4").Value = 24.06892892892927

# Row 5
$ws.Range("B5").Value = 0
$ws.Range("F5").Value = 24.19000000000034
$ws.Range("H5").Value = 0.03149955749924604
$ws.Range("I5").Value = 0.03149955749924604
$ws.Range("L5").Value = 25.80648916209751
$ws.Range("M5").Value = "[1.7161353014421437, 49.89684302275287]"
$ws.Range("N5").Value = 0.03633621323238789
$ws.Range("O5").Value = 0.03633621323238789
$ws.Range("P5").Value = 1.276763380738194
$ws.Range("Q5").Value = "[-0.06918422260157797, 2.6227109840779654]"
$ws.Range("R5").Value = 0.0624434432501102
$ws.Range("S5").Value = 0.0624434432501102
$ws.Range("T5").Value = 57.82651967175397
$ws.Range("U5").Value = "[44.35890810326862, 71.29413124023932]"
$ws.Range("V5").Value = [double]"4.0061731709784e-11"
$ws.Range("W5").Value = [double]"4.0061731709784e-11"
$ws.Range("X5").Value = 19.27451451451479
$ws.Range("Y5").Value = 14.09267267267287
$ws.Range("Z5").Value = 24.4563563563567

# Row 6
$ws.Range("F6").Value = 24.19000000000034
$ws.Range("H6").Value = 0.000119678506494636
$ws.Range("I6").Value = 0.000119678506494636
$ws.Range("L6").Value = 50.36362412809703
$ws.Range("M6").Value = "[25.003897088834762, 75.7233511673593]"
$ws.Range("N6").Value = 0.0002333088500574565
$ws.Range("O6").Value = 0.0002333088500574565
$ws.Range("P6").Value = 1.465447624197041
$ws.Range("Q6").Value = "[0.836500146000887, 2.094395102393195]"
$ws.Range("R6").Value = [double]"2.538723835798784e-05"
$ws.Range("S6").Value = [double]"2.538723835798784e-05"
$ws.Range("T6").Value = 64.17931835958022
$ws.Range("U6").Value = "[48.91601236944797, 79.44262434971247]"
$ws.Range("V6").Value = [double]"7.233280641116835e-11"
$ws.Range("W6").Value = [double]"7.233280641116835e-11"
$ws.Range("X6").Value = 18.54808808808835
$ws.Range("Y6").Value = 16.1266666666669
$ws.Range("Z6").Value = 20.9695095095098

# Row 7
$ws.Range("F7").Value = 24.19000000000034
$ws.Range("H7").Value = [double]"5.26713985891103e-05"
$ws.Range("I7").Value = [double]"5.26713985891103e-05"
$ws.Range("L7").Value = 43.75542970154144
$ws.Range("M7").Value = "[19.712403332825687, 67.7984560702572]"
$ws.Range("N7").Value = 0.000649111073027564
$ws.Range("O7").Value = 0.000649111073027564
$ws.Range("P7").Value = 0.723289599925578
$ws.Range("Q7").Value = "[0.1823947686768852, 1.2641844311742707]"
$ws.Range("R7").Value = 0.009905383267588386
$ws.Range("S7").Value = 0.009905383267588386
$ws.Range("T7").Value = 54.31985748598659
$ws.Range("U7").Value = "[41.80987440559885, 66.82984056637433]"
$ws.Range("V7").Value = [double]"2.909117391425298e-11"
$ws.Range("W7").Value = [double]"2.909117391425298e-11"
$ws.Range("X7").Value = 21.40536536536567
$ws.Range("Y7").Value = 19.32294294294322
$ws.Range("Z7").Value = 23.48778778778812

# Row 8
$ws.Range("F8").Value = 24.19000000000034
$ws.Range("H8").Value = 0.0009378884153542488
$ws.Range("I8").Value = 0.0009378884153542488
$ws.Range("L8").Value = 43.41940524967757
$ws.Range("M8").Value = "[15.731271597367268, 71.10753890198788]"
$ws.Range("N8").Value = 0.002832374208390664
$ws.Range("O8").Value = 0.002832374208390664
$ws.Range("P8").Value = 1.327079178993887
$ws.Range("Q8").Value = "[0.5975001042863468, 2.0566582537014266]"
$ws.Range("R8").Value = 0.0006527120661761199
$ws.Range("S8").Value = 0.0006527120661761199
$ws.Range("T8").Value = 67.55025626540296
$ws.Range("U8").Value = "[52.21900513541786, 82.88150739538807]"
$ws.Range("V8").Value = [double]"1.908873059619509e-11"
$ws.Range("W8").Value = [double]"1.908873059619509e-11"
$ws.Range("X8").Value = 19.08080080080107
$ws.Range("Y8").Value = 16.27195195195218
$ws.Range("Z8").Value = 21.88964964964996

# Row 9
$ws.Range("F9").Value = 24.19000000000034
$ws.Range("H9").Value = 0.006186995510041426
$ws.Range("I9").Value = 0.006186995510041426
$ws.Range("L9").Value = 33.52081059795327
$ws.Range("M9").Value = "[5.923423844703137, 61.118197351203406]"
$ws.Range("N9").Value = 0.01840084293680833
$ws.Range("O9").Value = 0.01840084293680833
$ws.Range("P9").Value = 0.9245527929483472
$ws.Range("Q9").Value = "[0.1823947686768843, 1.6667108172198102]"
$ws.Range("R9").Value = 0.015773177139772
$ws.Range("S9").Value = 0.015773177139772
$ws.Range("T9").Value = 66.0752555103874
$ws.Range("U9").Value = "[51.956404199339545, 80.19410682143524]"
$ws.Range("V9").Value = [double]"3.205880005907602e-12"
$ws.Range("W9").Value = [double]"3.205880005907602e-12"
$ws.Range("X9").Value = 20.6305105105108
$ws.Range("Y9").Value = 17.77323323323348
$ws.Range("Z9").Value = 23.48778778778812

# Row 10
$ws.Range("F10").Value = 24.19000000000034
$ws.Range("H10").Value = 0.00103465543651049
$ws.Range("I10").Value = 0.00103465543651049
$ws.Range("L10").Value = 39.34339993959943
$ws.Range("M10").Value = "[16.021188348439914, 62.665611530758945]"
$ws.Range("N10").Value = 0.00143137975866936
$ws.Range("O10").Value = 0.00143137975866936
$ws.Range("P10").Value = 1.440289725069194
$ws.Range("Q10").Value = "[0.685552751233808, 2.19502669890458]"
$ws.Range("R10").Value = 0.0003781036290455564
$ws.Range("S10").Value = 0.0003781036290455564
$ws.Range("T10").Value = 52.75551087744336
$ws.Range("U10").Value = "[38.749068109988144, 66.76195364489858]"
$ws.Range("V10").Value = [double]"1.388865911167159e-09"
$ws.Range("W10").Value = [double]"1.388865911167159e-09"
$ws.Range("X10").Value = 18.64494494494521
$ws.Range("Y10").Value = 15.73923923923946
$ws.Range("Z10").Value = 21.55065065065096

# Row 11
$ws.Range("F11").Value = 24.83000000000044
$ws.Range("H11").Value = 0.006447891480168666
$ws.Range("I11").Value = 0.006447891480168666
$ws.Range("L11").Value = 36.44011665385295
$ws.Range("M11").Value = "[10.875110394446416, 62.00512291325949]"
$ws.Range("N11").Value = 0.006218523301353418
$ws.Range("O11").Value = 0.006218523301353418
$ws.Range("P11").Value = 0.01886842434588498
$ws.Range("Q11").Value = "[-0.9434212172942322, 0.9811580659860022]"
$ws.Range("R11").Value = 0.9686726462886743
$ws.Range("S11").Value = 0.9686726462886743
$ws.Range("T11").Value = 55.01469578459257
$ws.Range("U11").Value = "[39.854882172498066, 70.17450939668707]"
$ws.Range("V11").Value = [double]"3.553661587218926e-09"
$ws.Range("W11").Value = [double]"3.553661587218926e-09"
$ws.Range("X11").Value = 24.75543543543587
$ws.Range("Y11").Value = 20.95264264264301
$ws.Range("Z11").Value = 28.55822822822874

# Row 12
$ws.Range("B12").Value = 0
$ws.Range("F12").Value = 24.83000000000044
$ws.Range("H12").Value = [double]"8.364243119451942e-06"
$ws.Range("I12").Value = [double]"8.364243119451942e-06"
$ws.Range("L12").Value = 59.00108800885423
$ws.Range("M12").Value = "[32.567207024445395, 85.43496899326306]"
$ws.Range("N12").Value = [double]"4.829451850918431e-05"
$ws.Range("O12").Value = [double]"4.829451850918431e-05"
$ws.Range("P12").Value = 0.2830263651882694
$ws.Range("Q12").Value = "[-0.2327105669325782, 0.798763297309117]"
$ws.Range("R12").Value = 0.2749039885462849
$ws.Range("S12").Value = 0.2749039885462849
$ws.Range("T12").Value = 64.16944662970516
$ws.Range("U12").Value = "[49.253940053636825, 79.0849532057735]"
$ws.Range("V12").Value = [double]"3.788125368942019e-11"
$ws.Range("W12").Value = [double]"3.788125368942019e-11"
$ws.Range("X12").Value = 23.71153153153195
$ws.Range("Y12").Value = 21.67343343343382
$ws.Range("Z12").Value = 25.74962962963009

# Row 13
$ws.Range("F13").Value = 24.83000000000044
$ws.Range("H13").Value = 0.6867260805483093
$ws.Range("I13").Value = 0.6867260805483093
$ws.Range("L13").Value = 9.917142030493757
$ws.Range("M13").Value = "[-20.310973454547792, 40.145257515535306]"
$ws.Range("N13").Value = 0.5121218984555305
$ws.Range("O13").Value = 0.5121218984555305
$ws.Range("P13").Value = -0.1257894956392311
$ws.Range("Q13").Value = "[-3.2202110883643122, 2.96863209708585]"
$ws.Range("R13").Value = 0.9351099027366394
$ws.Range("S13").Value = 0.9351099027366394
$ws.Range("T13").Value = 71.86998969886594
$ws.Range("U13").Value = "[55.8908285273056, 87.84915087042627]"
$ws.Range("V13").Value = [double]"1.04647401855118e-11"
$ws.Range("W13").Value = [double]"1.04647401855118e-11"
$ws.Range("X13").Value = 0.4970970970971074
$ws.Range("Y13").Value = -11.7314914914917
$ws.Range("Z13").Value = 12.72568568568591

# Row 14
$ws.Range("F14").Value = 24.83000000000044
$ws.Range("H14").Value = 0.02876368894042569
$ws.Range("I14").Value = 0.02876368894042569
$ws.Range("L14").Value = 27.91036184600318
$ws.Range("M14").Value = "[3.2391394391221553, 52.58158425288421]"
$ws.Range("N14").Value = 0.02749159025288583
$ws.Range("O14").Value = 0.02749159025288583
$ws.Range("P14").Value = 0.2704474156243464
$ws.Range("Q14").Value = "[-1.0188949146777704, 1.5597897459264631]"
$ws.Range("R14").Value = 0.6746935253782336
$ws.Range("S14").Value = 0.6746935253782336
$ws.Range("T14").Value = 54.02167834009504
$ws.Range("U14").Value = "[39.936949265520084, 68.10640741466999]"
$ws.Range("V14").Value = [double]"8.687768282555908e-10"
$ws.Range("W14").Value = [double]"8.687768282555908e-10"
$ws.Range("X14").Value = 23.76124124124167
$ws.Range("Y14").Value = 18.66599599599633
$ws.Range("Z14").Value = 28.856486486487

# Row 15
$ws.Range("F15").Value = 24.83000000000044
$ws.Range("H15").Value = 0.001015497445051583
$ws.Range("I15").Value = 0.001015497445051583
$ws.Range("L15").Value = 33.20724268198437
$ws.Range("M15").Value = "[11.906714253085227, 54.50777111088351]"
$ws.Range("N15").Value = 0.002982580554321101
$ws.Range("O15").Value = 0.002982580554321101
$ws.Range("P15").Value = -0.4905790329930007
$ws.Range("Q15").Value = "[-1.1572633598809245, 0.17610529389492324]"
$ws.Range("R15").Value = 0.1452895059372719
$ws.Range("S15").Value = 0.1452895059372719
$ws.Range("T15").Value = 50.15954207406931
$ws.Range("U15").Value = "[38.44925883087922, 61.86982531725941]"
$ws.Range("V15").Value = [double]"4.291211830320663e-11"
$ws.Range("W15").Value = [double]"4.291211830320663e-11"
$ws.Range("X15").Value = 1.938678678678713
$ws.Range("Y15").Value = -0.6959359359359478
$ws.Range("Z15").Value = 4.573293293293373

# Row 16
$ws.Range("F16").Value = 24.83000000000044
$ws.Range("H16").Value = [double]"2.971602498957893e-06"
$ws.Range("I16").Value = [double]"2.971602498957893e-06"
$ws.Range("L16").Value = 44.12626929151763
$ws.Range("M16").Value = "[23.20465401616423, 65.04788456687103]"
$ws.Range("N16").Value = 0.0001068989327797087
$ws.Range("O16").Value = 0.0001068989327797087
$ws.Range("P16").Value = -0.8931054190385401
$ws.Range("Q16").Value = "[-1.3711055024676178, -0.4151053356094625]"
$ws.Range("R16").Value = 0.0004831617890768136
$ws.Range("S16").Value = 0.0004831617890768136
$ws.Range("T16").Value = 57.53691534545195
$ws.Range("U16").Value = "[46.69462378287908, 68.37920690802483]"
$ws.Range("V16").Value = [double]"6.217248937900877e-14"
$ws.Range("W16").Value = [double]"6.217248937900877e-14"
$ws.Range("X16").Value = 3.529389389389454
$ws.Range("Y16").Value = 1.640420420420451
$ws.Range("Z16").Value = 5.418358358358457
